$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" values in column E (rows 16-22) are being re-arranged
# (previous periods removed / new periods added -> net effect is the list is
# reversed), and the corresponding due-date values in column F for the first
# and last row are swapped to match.

$ws.Range("E16").Value = "2311"
$ws.Range("E17").Value = "2310"
$ws.Range("E18").Value = "2309"
$ws.Range("E19").Value = "2308"
$ws.Range("E20").Value = "2307"
$ws.Range("E21").Value = "2306"
$ws.Range("E22").Value = "2305"

$ws.Range("F16").Value = 43307
$ws.Range("F22").Value = 44028
